# The deck currently carries two leftover theme parts from an earlier
# template switch: the slide master / design ("Integral") and a second,
# otherwise-unused theme ("Office Theme") that only the notes master
# still points at. This commit swaps the two designs around: the
# slides should render with the plain "Office Theme" palette, while
# the (COM-unreachable) notes-only theme keeps the old "Integral"
# palette. The only thing PowerPoint's object model exposes for
# editing a design's theme in place is its ThemeColorScheme, so the
# swap is performed colour-by-colour on the presentation's one (and
# only) reachable Design/SlideMaster theme.

function Make-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The Design backing the slide master (and therefore every slide) —
# this is the theme object PowerPoint's COM automation can actually
# reach and persist edits to.
$design = $p.Designs.Item(1)
$theme  = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Target palette: the stock Office theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in ThemeColorScheme's fixed 1-12 order.
$colors.Item(1).RGB  = Make-RGB 0x00 0x00 0x00   # dk1       #000000
$colors.Item(2).RGB  = Make-RGB 0xFF 0xFF 0xFF   # lt1       #FFFFFF
$colors.Item(3).RGB  = Make-RGB 0x44 0x54 0x6A   # dk2       #44546A
$colors.Item(4).RGB  = Make-RGB 0xE7 0xE6 0xE6   # lt2       #E7E6E6
$colors.Item(5).RGB  = Make-RGB 0x5B 0x9B 0xD5   # accent1   #5B9BD5
$colors.Item(6).RGB  = Make-RGB 0xED 0x7D 0x31   # accent2   #ED7D31
$colors.Item(7).RGB  = Make-RGB 0xA5 0xA5 0xA5   # accent3   #A5A5A5
$colors.Item(8).RGB  = Make-RGB 0xFF 0xC0 0x00   # accent4   #FFC000
$colors.Item(9).RGB  = Make-RGB 0x44 0x72 0xC4   # accent5   #4472C4
$colors.Item(10).RGB = Make-RGB 0x70 0xAD 0x47   # accent6   #70AD47
$colors.Item(11).RGB = Make-RGB 0x05 0x63 0xC1   # hlink     #0563C1
$colors.Item(12).RGB = Make-RGB 0x95 0x4F 0x72   # folHlink  #954F72

# Best-effort: some hosts also let the design/theme display name be
# updated directly. Harmless no-op where the property is read-only.
try { $design.Name = "Office Theme" } catch {}
try { $theme.Name  = "Office Theme" } catch {}
